# Update countries & provincias Spain
# - Re-sort the "Angola" entry to sit right after "Mongolia" in the country
#   list (this shifts the countries that used to sit between "Liberia" and
#   "Angola" down by one row), and refresh the day's COVID-19 case figures.
# - Bump the "Datos actualizados..." timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 00:52"

# --- Updated totals for a few already-listed countries ---------------------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 334745
$ws.Range("C4").Value = 23388
$ws.Range("E4").Value = 307931
$ws.Range("G4").Value = 1121
$ws.Range("H4").Value = 9572

# Row 24
$ws.Range("E24").Value = 3335
$ws.Range("G24").Value = 7
$ws.Range("H24").Value = 37

# Row 27
$ws.Range("B27").Value = 4587
$ws.Range("C27").Value = 115
$ws.Range("E27").Value = 4424

# --- Country list re-sort: "Angola" now ranks right after "Mongolia" -------
# Rows 173-184 shift to the country that used to be one row below, each one
# carrying that country's current figures.
$rows = @(173, 174, 175, 176, 177, 178, 179, 180, 181, 182, 183, 184)
$countryNames = @{
    173 = "Angola"
    174 = "Liberia"
    175 = "Granada"
    176 = "Fiyi"
    177 = "Sudan"
    178 = "Laos"
    179 = "Groenlandia"
    180 = "Curazao"
    181 = "San Cristobal y Nieves"
    182 = "Seychelles"
    183 = "Mozambique"
    184 = "Surinam"
}
$countryStats = @{
    173 = @{ B=14; C=4;  D=2; E=10; F=0; G=0; H=2 }
    174 = @{ B=13; C=3;  D=3; E=7;  F=0; G=2; H=3 }
    175 = @{ B=12; C=0;  D=0; E=12; F=2; G=0; H=0 }
    176 = @{ B=12; C=0;  D=0; E=12; F=0; G=0; H=0 }
    177 = @{ B=12; C=2;  D=2; E=8;  F=0; G=0; H=2 }
    178 = @{ B=11; C=1;  D=0; E=11; F=0; G=0; H=0 }
    179 = @{ B=11; C=0;  D=3; E=8;  F=0; G=0; H=0 }
    180 = @{ B=11; C=0;  D=5; E=5;  F=0; G=0; H=1 }
    181 = @{ B=10; C=1;  D=0; E=10; F=0; G=0; H=0 }
    182 = @{ B=10; C=0;  D=0; E=10; F=0; G=0; H=0 }
    183 = @{ B=10; C=0;  D=1; E=9;  F=0; G=0; H=0 }
    184 = @{ B=10; C=0;  D=0; E=9;  F=0; G=0; H=1 }
}

foreach ($r in $rows) {
    $ws.Range("A$r").Value = $countryNames[$r]
    $stats = $countryStats[$r]
    foreach ($col in @("B", "C", "D", "E", "F", "G", "H")) {
        $ws.Range("$col$r").Value = $stats[$col]
    }
}
